# Append a new log row (row 8) to the "Proximity" sheet, matching the
# existing rows: a Bedroom Door ENTER event.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proximity")

# Column A holds an ISO-like date string ("2026-02-01") which Excel's
# auto-detection would otherwise convert into a real date serial number.
# Briefly force the cell to Text format before writing the value (like
# pre-formatting a cell as Text in the UI), then clear the formatting
# again so the cell ends up with no special style applied - exactly like
# the sibling cells A2:A7 above it.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2026-02-01"
$ws.Range("A8").ClearFormats()

$ws.Range("B8").Value = "15:07:48"
$ws.Range("C8").Value = "15:00"
$ws.Range("D8").Value = "Bedroom Door"
$ws.Range("E8").Value = "ENTER"
$ws.Range("F8").Value = "User ENTERED Bedroom"
